$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Fix the "ccm_ui_remove" -> "cmm_ui_remove" typo in row 4 ---
$ws.Range("A4").Value = "cmm_ui_remove"

# --- 2) New row 5 ("Unknown Mod" hint) re-uses the big, special JP font that used to ---
#        live on C4 ("無し"). Copy that formatting onto C5 BEFORE we touch C4, so the
#        style carries over to the new row intact.
$ws.Range("C4").Copy() | Out-Null
$ws.Range("C5").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# --- 3) Give the rest of row 5 the same plain formatting as the rest of row 4 ---
$ws.Range("A4").Copy() | Out-Null
$ws.Range("A5").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("B4").Copy() | Out-Null
$ws.Range("B5").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("D4").Copy() | Out-Null
$ws.Range("D5").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# --- 4) C4 ("無し") drops the special big font and becomes plain, like D4 ---
$ws.Range("D4").Copy() | Out-Null
$ws.Range("C4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# --- 5) Fill in the new "Unknown Mod" row ---
$ws.Range("A5").Value = "cmm_ui_unknown"
$ws.Range("D5").Value = "Unknown Mod"

# C5 gets a two-run rich string: "不明な" in the normal (inherited) font, "MOD" forced
# to Cascadia Code 16pt (matching the look of the other Cascadia-Code cells).
$ws.Range("C5").Value = "不明なMOD"
$modChars = $ws.Range("C5").Characters(4, 3)
$modChars.Font.Name = "Cascadia Code"
$modChars.Font.Size = 16

# --- 6) Row 5 should be the same (auto/"big font") height as rows 3-4 ---
$ws.Range("A5").RowHeight = 23.25

# --- 7) The special JP font used for these callouts is renamed (was 微软雅黑) ---
$ws.Range("C5").Font.Name = "宋体"

# --- 8) Matches the author's final cursor position ---
$ws.Range("D9").Select() | Out-Null
